$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(74, 8).Value = 4550.7144  # H74: 4431.875 -> 4550.7144
$ws.Cells.Item(74, 9).Value = 4691.6  # I74: 4509.6665 -> 4691.6
$ws.Cells.Item(74, 11).Value = 4691.6  # K74: 4509.6665 -> 4691.6
$ws.Cells.Item(74, 13).Value = -3755.6  # M74: -3573.6665 -> -3755.6

$ws.Cells.Item(77, 8).Value = 4550.7144  # H77: 4431.875 -> 4550.7144
$ws.Cells.Item(77, 9).Value = 4691.6  # I77: 4509.6665 -> 4691.6
$ws.Cells.Item(77, 11).Value = 23458  # K77: 22548.3325 -> 23458
$ws.Cells.Item(77, 13).Value = -18778  # M77: -17868.3325 -> -18778

$ws.Cells.Item(86, 8).Value = 3833.3333  # H86: 5250 -> 3833.3333
$ws.Cells.Item(86, 9).Value = 2250  # I86: 3500 -> 2250
$ws.Cells.Item(86, 11).Value = 2250  # K86: 3500 -> 2250
$ws.Cells.Item(86, 13).Value = -1127  # M86: -2377 -> -1127

$ws.Cells.Item(89, 8).Value = 3833.3333  # H89: 5250 -> 3833.3333
$ws.Cells.Item(89, 9).Value = 2250  # I89: 3500 -> 2250
$ws.Cells.Item(89, 11).Value = 11250  # K89: 17500 -> 11250
$ws.Cells.Item(89, 13).Value = -5634  # M89: -11884 -> -5634

$ws.Cells.Item(121, 8).Value = 2024.5  # H121: 1382.7693 -> 2024.5
$ws.Cells.Item(121, 10).Value = 2024.5  # J121: 1382.7693 -> 2024.5
$ws.Cells.Item(121, 12).Value = 6073.5  # L121: 4148.3079 -> 6073.5
$ws.Cells.Item(121, 14).Value = -9567.5  # N121: -7642.3079 -> -9567.5

$ws.Cells.Item(137, 8).Value = 2142.4443  # H137: 2118.6553 -> 2142.4443
$ws.Cells.Item(137, 9).Value = 1996.4  # I137: 1974.5 -> 1996.4
$ws.Cells.Item(137, 10).Value = 2325  # J137: 2354.5454 -> 2325
$ws.Cells.Item(137, 11).Value = 5989.200000000001  # K137: 5923.5 -> 5989.200000000001
$ws.Cells.Item(137, 12).Value = 6975  # L137: 7063.6362 -> 6975
$ws.Cells.Item(137, 13).Value = -3439.200000000001  # M137: -3373.5 -> -3439.200000000001
$ws.Cells.Item(137, 14).Value = -12075  # N137: -12163.6362 -> -12075

$ws.Cells.Item(138, 8).Value = 1950.5428  # H138: 2013.9487 -> 1950.5428
$ws.Cells.Item(138, 9).Value = 1317.1818  # I138: 1372.5238 -> 1317.1818
$ws.Cells.Item(138, 10).Value = 3022.3845  # J138: 2762.2778 -> 3022.3845
$ws.Cells.Item(138, 11).Value = 3951.5454  # K138: 4117.5714 -> 3951.5454
$ws.Cells.Item(138, 12).Value = 9067.1535  # L138: 8286.8334 -> 9067.1535
$ws.Cells.Item(138, 13).Value = 1188.4546  # M138: 1022.4286 -> 1188.4546
$ws.Cells.Item(138, 14).Value = -19347.1535  # N138: -18566.8334 -> -19347.1535

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 3236.7903  # H32: 3407 -> 3236.7903
$ws.Cells.Item(32, 9).Value = 2332.7083  # I32: 2474.8865 -> 2332.7083
$ws.Cells.Item(32, 11).Value = 2332.7083  # K32: 2474.8865 -> 2332.7083
$ws.Cells.Item(32, 13).Value = -2045.7083  # M32: -2187.8865 -> -2045.7083

$ws.Cells.Item(74, 8).Value = 1631.1  # H74: 1680.1111 -> 1631.1
$ws.Cells.Item(74, 9).Value = 1535.8889  # I74: 1579.125 -> 1535.8889
$ws.Cells.Item(74, 11).Value = 1535.8889  # K74: 1579.125 -> 1535.8889
$ws.Cells.Item(74, 13).Value = -661.8888999999999  # M74: -705.125 -> -661.8888999999999

$ws.Cells.Item(77, 8).Value = 1631.1  # H77: 1680.1111 -> 1631.1
$ws.Cells.Item(77, 9).Value = 1535.8889  # I77: 1579.125 -> 1535.8889
$ws.Cells.Item(77, 11).Value = 7679.4445  # K77: 7895.625 -> 7679.4445
$ws.Cells.Item(77, 13).Value = -3311.4445  # M77: -3527.625 -> -3311.4445

$ws.Cells.Item(122, 8).Value = 1701.5  # H122: 1725.8572 -> 1701.5
$ws.Cells.Item(122, 9).Value = 1725.381  # I122: 1752.15 -> 1725.381
$ws.Cells.Item(122, 11).Value = 5176.143  # K122: 5256.450000000001 -> 5176.143
$ws.Cells.Item(122, 13).Value = -2726.143  # M122: -2806.450000000001 -> -2726.143

$ws.Cells.Item(132, 8).Value = 4113.4287  # H132: 4149 -> 4113.4287
$ws.Cells.Item(132, 9).Value = 3949.25  # I132: 3965.6667 -> 3949.25
$ws.Cells.Item(132, 11).Value = 11847.75  # K132: 11897.0001 -> 11847.75
$ws.Cells.Item(132, 13).Value = -9317.75  # M132: -9367.000100000001 -> -9317.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(99, 8).Value = 1357.5714  # H99: 1234.625 -> 1357.5714
$ws.Cells.Item(99, 9).Value = 1278.6666  # I99: 1121.3334 -> 1278.6666
$ws.Cells.Item(99, 10).Value = 1499.6  # J99: 1574.5 -> 1499.6
$ws.Cells.Item(99, 11).Value = 1278.6666  # K99: 1121.3334 -> 1278.6666
$ws.Cells.Item(99, 12).Value = 1499.6  # L99: 1574.5 -> 1499.6
$ws.Cells.Item(99, 13).Value = 219.3334  # M99: 376.6666 -> 219.3334
$ws.Cells.Item(99, 14).Value = -4495.6  # N99: -4570.5 -> -4495.6

$ws.Cells.Item(105, 8).Value = 1947  # H105: 2009.16 -> 1947
$ws.Cells.Item(105, 9).Value = 1953.6666  # I105: 2090.3684 -> 1953.6666
$ws.Cells.Item(105, 10).Value = 1919  # J105: 1752 -> 1919
$ws.Cells.Item(105, 11).Value = 1953.6666  # K105: 2090.3684 -> 1953.6666
$ws.Cells.Item(105, 12).Value = 1919  # L105: 1752 -> 1919
$ws.Cells.Item(105, 13).Value = -206.6666  # M105: -343.3683999999998 -> -206.6666
$ws.Cells.Item(105, 14).Value = -5413  # N105: -5246 -> -5413

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 4160.6665  # H31: 4701.077 -> 4160.6665
$ws.Cells.Item(31, 9).Value = 1827.875  # I31: 1806 -> 1827.875
$ws.Cells.Item(31, 10).Value = 6026.9  # J31: 5987.778 -> 6026.9
$ws.Cells.Item(31, 11).Value = 1827.875  # K31: 1806 -> 1827.875
$ws.Cells.Item(31, 12).Value = 6026.9  # L31: 5987.778 -> 6026.9
$ws.Cells.Item(31, 13).Value = -1532.875  # M31: -1511 -> -1532.875
$ws.Cells.Item(31, 14).Value = -6616.9  # N31: -6577.778 -> -6616.9

$ws.Cells.Item(34, 8).Value = 4160.6665  # H34: 4701.077 -> 4160.6665
$ws.Cells.Item(34, 9).Value = 1827.875  # I34: 1806 -> 1827.875
$ws.Cells.Item(34, 10).Value = 6026.9  # J34: 5987.778 -> 6026.9
$ws.Cells.Item(34, 11).Value = 1827.875  # K34: 1806 -> 1827.875
$ws.Cells.Item(34, 12).Value = 6026.9  # L34: 5987.778 -> 6026.9
$ws.Cells.Item(34, 13).Value = -1625.875  # M34: -1604 -> -1625.875
$ws.Cells.Item(34, 14).Value = -6430.9  # N34: -6391.778 -> -6430.9

$ws.Cells.Item(94, 8).Value = 1187.0714  # H94: 1285.4166 -> 1187.0714
$ws.Cells.Item(94, 9).Value = 1268.8334  # I94: 1322.6 -> 1268.8334
$ws.Cells.Item(94, 10).Value = 1125.75  # J94: 1258.8572 -> 1125.75
$ws.Cells.Item(94, 11).Value = 1268.8334  # K94: 1322.6 -> 1268.8334
$ws.Cells.Item(94, 12).Value = 1125.75  # L94: 1258.8572 -> 1125.75
$ws.Cells.Item(94, 13).Value = -817.8334  # M94: -871.5999999999999 -> -817.8334
$ws.Cells.Item(94, 14).Value = -2027.75  # N94: -2160.8572 -> -2027.75

$ws.Cells.Item(99, 8).Value = 2099.7778  # H99: 1949.6666 -> 2099.7778
$ws.Cells.Item(99, 9).Value = 1899.6666  # I99: 1839.6 -> 1899.6666
$ws.Cells.Item(99, 11).Value = 1899.6666  # K99: 1839.6 -> 1899.6666
$ws.Cells.Item(99, 13).Value = -401.6666  # M99: -341.5999999999999 -> -401.6666

$ws.Cells.Item(122, 8).Value = 2454.923  # H122: 1867.8572 -> 2454.923
$ws.Cells.Item(122, 9).Value = 2800.4  # I122: 1959.1111 -> 2800.4
$ws.Cells.Item(122, 10).Value = 2239  # J122: 1799.4166 -> 2239
$ws.Cells.Item(122, 11).Value = 8401.200000000001  # K122: 5877.3333 -> 8401.200000000001
$ws.Cells.Item(122, 12).Value = 6717  # L122: 5398.2498 -> 6717
$ws.Cells.Item(122, 13).Value = -5951.200000000001  # M122: -3427.3333 -> -5951.200000000001
$ws.Cells.Item(122, 14).Value = -11617  # N122: -10298.2498 -> -11617

$ws.Cells.Item(126, 8).Value = 2099.7778  # H126: 1949.6666 -> 2099.7778
$ws.Cells.Item(126, 9).Value = 1899.6666  # I126: 1839.6 -> 1899.6666
$ws.Cells.Item(126, 11).Value = 5698.9998  # K126: 5518.799999999999 -> 5698.9998
$ws.Cells.Item(126, 13).Value = -3228.9998  # M126: -3048.799999999999 -> -3228.9998

$ws.Cells.Item(132, 8).Value = 2222.8462  # H132: 2446.7273 -> 2222.8462
$ws.Cells.Item(132, 9).Value = 1065.6  # I132: 1092.5454 -> 1065.6
$ws.Cells.Item(132, 11).Value = 3196.8  # K132: 3277.6362 -> 3196.8
$ws.Cells.Item(132, 13).Value = -666.7999999999997  # M132: -747.6361999999999 -> -666.7999999999997

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(98, 8).Value = 554  # H98: 522.7143 -> 554
$ws.Cells.Item(98, 9).Value = 0  # I98: 400 -> 0
$ws.Cells.Item(98, 10).Value = 554  # J98: 543.1667 -> 554
$ws.Cells.Item(98, 11).Value = 0  # K98: 1200 -> 0
$ws.Cells.Item(98, 12).Value = 1662  # L98: 1629.5001 -> 1662
$ws.Cells.Item(98, 13).ClearContents()  # M98: remove (was 298)
$ws.Cells.Item(98, 14).Value = -4658  # N98: -4625.5001 -> -4658

$ws.Cells.Item(122, 8).Value = 1213.9166  # H122: 1278.909 -> 1213.9166
$ws.Cells.Item(122, 9).Value = 947  # I122: 1016.8 -> 947
$ws.Cells.Item(122, 10).Value = 1480.8334  # J122: 1497.3334 -> 1480.8334
$ws.Cells.Item(122, 11).Value = 8523  # K122: 9151.199999999999 -> 8523
$ws.Cells.Item(122, 12).Value = 13327.5006  # L122: 13476.0006 -> 13327.5006
$ws.Cells.Item(122, 13).Value = -6073  # M122: -6701.199999999999 -> -6073
$ws.Cells.Item(122, 14).Value = -18227.5006  # N122: -18376.0006 -> -18227.5006

$ws.Cells.Item(129, 8).Value = 32411.957  # H129: 33874.227 -> 32411.957
$ws.Cells.Item(129, 9).Value = 532.2857  # I129: 599.2 -> 532.2857
$ws.Cells.Item(129, 10).Value = 46359.312  # J129: 43661 -> 46359.312
$ws.Cells.Item(129, 11).Value = 1596.8571  # K129: 1797.6 -> 1596.8571
$ws.Cells.Item(129, 12).Value = 139077.936  # L129: 130983 -> 139077.936
$ws.Cells.Item(129, 13).Value = 3403.1429  # M129: 3202.4 -> 3403.1429
$ws.Cells.Item(129, 14).Value = -149077.936  # N129: -140983 -> -149077.936

$ws.Cells.Item(131, 8).Value = 18545012  # H131: 15647479 -> 18545012
$ws.Cells.Item(131, 9).Value = 55556024  # I131: 62500424 -> 55556024
$ws.Cells.Item(131, 10).Value = 39505.445  # J131: 29829.916 -> 39505.445
$ws.Cells.Item(131, 11).Value = 166668072  # K131: 187501272 -> 166668072
$ws.Cells.Item(131, 12).Value = 118516.335  # L131: 89489.74800000001 -> 118516.335
$ws.Cells.Item(131, 13).Value = -166663032  # M131: -187496232 -> -166663032
$ws.Cells.Item(131, 14).Value = -128596.335  # N131: -99569.74800000001 -> -128596.335

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(92, 8).Value = 15998.5  # H92: 12091.625 -> 15998.5
$ws.Cells.Item(92, 10).Value = 15998.5  # J92: 12091.625 -> 15998.5
$ws.Cells.Item(92, 12).Value = 15998.5  # L92: 12091.625 -> 15998.5
$ws.Cells.Item(92, 14).Value = -19742.5  # N92: -15835.625 -> -19742.5

$ws.Cells.Item(102, 8).Value = 2066.9  # H102: 2122.4736 -> 2066.9
$ws.Cells.Item(102, 9).Value = 2439.1  # I102: 2611.2222 -> 2439.1
$ws.Cells.Item(102, 10).Value = 1694.7  # J102: 1682.6 -> 1694.7
$ws.Cells.Item(102, 11).Value = 2439.1  # K102: 2611.2222 -> 2439.1
$ws.Cells.Item(102, 12).Value = 1694.7  # L102: 1682.6 -> 1694.7
$ws.Cells.Item(102, 13).Value = -817.0999999999999  # M102: -989.2222000000002 -> -817.0999999999999
$ws.Cells.Item(102, 14).Value = -4938.7  # N102: -4926.6 -> -4938.7

$ws.Cells.Item(122, 8).Value = 2301.9167  # H122: 2481.6667 -> 2301.9167
$ws.Cells.Item(122, 9).Value = 2169.2222  # I122: 2267 -> 2169.2222
$ws.Cells.Item(122, 10).Value = 2700  # J122: 2750 -> 2700
$ws.Cells.Item(122, 11).Value = 6507.6666  # K122: 6801 -> 6507.6666
$ws.Cells.Item(122, 12).Value = 8100  # L122: 8250 -> 8100
$ws.Cells.Item(122, 13).Value = -4057.6666  # M122: -4351 -> -4057.6666
$ws.Cells.Item(122, 14).Value = -13000  # N122: -13150 -> -13000

$ws.Cells.Item(126, 8).Value = 38525.215  # H126: 39848.184 -> 38525.215
$ws.Cells.Item(126, 9).Value = 3023.6  # I126: 3034.3157 -> 3023.6
$ws.Cells.Item(126, 10).Value = 127279.25  # J126: 127281.125 -> 127279.25
$ws.Cells.Item(126, 11).Value = 9070.799999999999  # K126: 9102.947100000001 -> 9070.799999999999
$ws.Cells.Item(126, 12).Value = 381837.75  # L126: 381843.375 -> 381837.75
$ws.Cells.Item(126, 13).Value = -6600.799999999999  # M126: -6632.947100000001 -> -6600.799999999999
$ws.Cells.Item(126, 14).Value = -386777.75  # N126: -386783.375 -> -386777.75

$ws.Cells.Item(132, 8).Value = 3049.6086  # H132: 3291.0527 -> 3049.6086
$ws.Cells.Item(132, 9).Value = 2643.923  # I132: 2866 -> 2643.923
$ws.Cells.Item(132, 10).Value = 3577  # J132: 3763.3333 -> 3577
$ws.Cells.Item(132, 11).Value = 7931.768999999999  # K132: 8598 -> 7931.768999999999
$ws.Cells.Item(132, 12).Value = 10731  # L132: 11289.9999 -> 10731
$ws.Cells.Item(132, 13).Value = -5401.768999999999  # M132: -6068 -> -5401.768999999999
$ws.Cells.Item(132, 14).Value = -15791  # N132: -16349.9999 -> -15791

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(16, 8).Value = 31399.666  # H16: 8114.75 -> 31399.666
$ws.Cells.Item(16, 9).Value = 31399.666  # I16: 9202.571 -> 31399.666
$ws.Cells.Item(16, 10).Value = 0  # J16: 500 -> 0
$ws.Cells.Item(16, 11).Value = 31399.666  # K16: 9202.571 -> 31399.666
$ws.Cells.Item(16, 12).Value = 0  # L16: 500 -> 0
$ws.Cells.Item(16, 13).Value = -31229.666  # M16: -9032.571 -> -31229.666
$ws.Cells.Item(16, 14).ClearContents()  # N16: remove (was -840)

$ws.Cells.Item(46, 8).Value = 1559  # H46: 1540.8182 -> 1559
$ws.Cells.Item(46, 9).Value = 800  # I46: 600 -> 800
$ws.Cells.Item(46, 11).Value = 800  # K46: 600 -> 800
$ws.Cells.Item(46, 13).Value = -612  # M46: -412 -> -612

$ws.Cells.Item(55, 8).Value = 239.5  # H55: 248.18182 -> 239.5
$ws.Cells.Item(55, 9).Value = 192.8  # I55: 216.33333 -> 192.8
$ws.Cells.Item(55, 10).Value = 272.85715  # J55: 260.125 -> 272.85715
$ws.Cells.Item(55, 11).Value = 192.8  # K55: 216.33333 -> 192.8
$ws.Cells.Item(55, 12).Value = 272.85715  # L55: 260.125 -> 272.85715
$ws.Cells.Item(55, 13).Value = -19.80000000000001  # M55: -43.33332999999999 -> -19.80000000000001
$ws.Cells.Item(55, 14).Value = -618.85715  # N55: -606.125 -> -618.85715

$ws.Cells.Item(122, 8).Value = 7984.615  # H122: 8709.091 -> 7984.615
$ws.Cells.Item(122, 10).Value = 8400  # J122: 9500 -> 8400
$ws.Cells.Item(122, 12).Value = 25200  # L122: 28500 -> 25200
$ws.Cells.Item(122, 14).Value = -30100  # N122: -33400 -> -30100

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(122, 8).Value = 56945.715  # H122: 47130.883 -> 56945.715
$ws.Cells.Item(122, 9).Value = 72072.82000000001  # I122: 56913.285 -> 72072.82000000001
$ws.Cells.Item(122, 11).Value = 216218.46  # K122: 170739.855 -> 216218.46
$ws.Cells.Item(122, 13).Value = -213768.46  # M122: -168289.855 -> -213768.46

$ws.Cells.Item(136, 8).Value = 2380.9644  # H136: 2401.6333 -> 2380.9644
$ws.Cells.Item(136, 9).Value = 1509.1177  # I136: 1676.8 -> 1509.1177
$ws.Cells.Item(136, 10).Value = 3728.3635  # J136: 3851.3 -> 3728.3635
$ws.Cells.Item(136, 11).Value = 4527.3531  # K136: 5030.4 -> 4527.3531
$ws.Cells.Item(136, 12).Value = 11185.0905  # L136: 11553.9 -> 11185.0905
$ws.Cells.Item(136, 13).Value = -1977.3531  # M136: -2480.4 -> -1977.3531
$ws.Cells.Item(136, 14).Value = -16285.0905  # N136: -16653.9 -> -16285.0905
